# Auto-generated Excel COM-interop script
# Updates hashcode values (column B) in the "hashcode.csv" worksheet
# to match the new values recorded in the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B9").Value = "e93effb58e5970f605ae07ea0fd6480b"
$ws.Range("B11").Value = "27727890b1e1b6a6913a83c5b04b29a4"
$ws.Range("B17").Value = "bb451ec4926ef9a76c82b3a70560c0a5"
$ws.Range("B29").Value = "283a43421d3d619311cfd0592b2dd6eb"
$ws.Range("B44").Value = "64d97a2435ca528474a9ee1b62a5969d"
$ws.Range("B74").Value = "7ab7fef2fd4db72bbdb0720aafcbc718"
$ws.Range("B89").Value = "677808ed7f974be62cdfb69b4daed013"
$ws.Range("B99").Value = "7295799e6758bfbfe9f01c6adf1aca08"
$ws.Range("B110").Value = "1bd5e3b761a52acf1e20b0c69324b0d1"
$ws.Range("B121").Value = "a75a4d9af7223344b490c2aca5cdac25"
$ws.Range("B126").Value = "51bbf56d85cc17f3c8cb856bf4fd262d"
$ws.Range("B133").Value = "6fb7a167831c59c63f682576093f9892"
$ws.Range("B136").Value = "145f6cdd9e574970a49058607a4c57c6"
$ws.Range("B159").Value = "dbfc21f7e94c2499a7e91e097f364003"
$ws.Range("B161").Value = "10f1715cd7ab53d5a3f38c26ac1e512f"
$ws.Range("B168").Value = "b59d55c420b531bf2814747715b21456"
$ws.Range("B169").Value = "d8e2d3b430620fbcc36650018a5d213d"
$ws.Range("B191").Value = "3a13e7d435e81d1a9016877eee3af917"
$ws.Range("B246").Value = "05a60c2804215dc6c1e82593bfcca0ec"
$ws.Range("B276").Value = "50ca57d3bac2b87bc65ddb88545854d8"
$ws.Range("B278").Value = "6ca2b727497da9da297e10d0e74f11fc"
$ws.Range("B281").Value = "7f6ab24a2600337270ff3e0396ae3efd"
$ws.Range("B293").Value = "8cb4f938f3e6a3f50370cb099b1625d5"
$ws.Range("B302").Value = "0f1ef506e706195dbd93c49065f789b1"
$ws.Range("B339").Value = "4355b8ccd9f3d91560badc347230afcd"
$ws.Range("B345").Value = "1d0565d3900a06151050ed3f0730ef7c"
$ws.Range("B410").Value = "4c457074098f96716ec4fb9f40496aa5"
$ws.Range("B446").Value = "7c50a6ebad1e7320fab77d231a4aad3b"
$ws.Range("B460").Value = "ef3bb11c9a11290215fab20c3653025e"
$ws.Range("B480").Value = "54047bec7956934d2f51b05c58bf2b32"
$ws.Range("B500").Value = "90638a5840cb2ea45547ac598d99705e"
$ws.Range("B501").Value = "10add39a694426657601535a2ecb2c04"
$ws.Range("B515").Value = "20970741bb8f5220d99c759f67734917"
$ws.Range("B516").Value = "2f7c284edd9cf2e6476c05e9d963bf7e"
$ws.Range("B517").Value = "d58681c86cbed19c395aab18d70338ab"
$ws.Range("B550").Value = "8aab137630c87b0adee966d8555f7e13"
$ws.Range("B566").Value = "27c4e453f2e3cd5748c259466fa6b621"
$ws.Range("B572").Value = "0751fcd52a01e68b0dea88477cc78546"
$ws.Range("B616").Value = "078638d89707ef761041c1aa1f6eb798"
$ws.Range("B627").Value = "0225aa8685f6b6a513936ce0d53587e9"
$ws.Range("B629").Value = "8e135f17d024197e1fee484b3eb87bd1"
$ws.Range("B649").Value = "3e72d49f2d5a1c266973b510c1bc866a"
$ws.Range("B655").Value = "d6d55401dea2dc036bcb028447293785"
$ws.Range("B700").Value = "c1be0d083ce0ad19eb1f14e63dd5771f"
$ws.Range("B756").Value = "34bc1b8b8dd8361c7d36d04fc4d5341b"
$ws.Range("B761").Value = "becaa9f87c93d226e1846c5199e485c4"
$ws.Range("B768").Value = "fa3438559eb36bcd278952eeb9ffd616"
$ws.Range("B786").Value = "e04abe38bb57d5e3316669455587b93b"
$ws.Range("B816").Value = "e156ff61a68c1b859d559b0ba2bd01c0"
$ws.Range("B819").Value = "ddcecae74f700d34aeb688e4eafe9966"
$ws.Range("B825").Value = "74f20965bca711405d4b5008fd88b85c"
$ws.Range("B827").Value = "7c0d8b2c888ea89da57dac14fe891e28"
$ws.Range("B830").Value = "878f501c6fcfbb24100b756563e49341"
$ws.Range("B835").Value = "44a1dc031076aedec8ddf2465a2c79d5"
$ws.Range("B855").Value = "b36fc7e153b7129ca0d1943ba1e30c19"
$ws.Range("B862").Value = "56ad9242b497ae392e8130d0697a5abd"
$ws.Range("B869").Value = "3fd20f2c7999ccdb99f29dcaa7ffe607"
$ws.Range("B874").Value = "d878f735a89572d2273c1e98708e28dd"
$ws.Range("B904").Value = "61df70027f6718b3d5068d6166901e38"
$ws.Range("B928").Value = "b32936677370cb45562c08f26acd3ef0"

